$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cross Section Properties")
$ws.Activate()

# Row 6: the stringer input for E6 flips sign (moved to the other side of the section)
$ws.Range("E6").Value = -0.0625

# Row 11: the hardcoded geometry for this new/repositioned stringer is replaced by
# explicit coordinate formulas (negative values -> moved relative to the reference axis),
# and the old cell formatting (fill style index 5) is cleared back to Normal.
$ws.Range("D11").Style = "Normal"
$ws.Range("D11").Formula = "=-1.6046754518"
$ws.Range("E11").Style = "Normal"
$ws.Range("E11").Formula = "=-0.092597515"

# Move the active selection/view to match the saved workbook state
$ws.Range("E7").Select()
